$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 18: F 81 -> 82, H 118 -> 119
$ws.Range("F18").Value = 82
$ws.Range("H18").Value = 119

# Row 33: F 19 -> 20, H 31 -> 32
$ws.Range("F33").Value = 20
$ws.Range("H33").Value = 32

# Row 34: F 17 -> 18, H 20 -> 21
$ws.Range("F34").Value = 18
$ws.Range("H34").Value = 21

# Row 36: F 76 -> 78, H 108 -> 110
$ws.Range("F36").Value = 78
$ws.Range("H36").Value = 110

# Row 43: F 26 -> 27, H 29 -> 30
$ws.Range("F43").Value = 27
$ws.Range("H43").Value = 30

# Row 58: F 5 -> 6, H 5 -> 6
$ws.Range("F58").Value = 6
$ws.Range("H58").Value = 6
